# Refresh the cryptocurrency price/volume snapshot (GitHub Actions cron
# update). Price/Volume(1h) cells are stored as plain text in this sheet
# (e.g. "66.960.43", "  -2.01%  ") rather than numbers, so most updates are
# simple literal string writes. A handful of new Price values are fully
# numeric-looking text (e.g. "7.70", "0.168") that Excel's COM layer would
# otherwise auto-coerce into a real number (dropping the trailing zero /
# losing the original text formatting) - those are written with a leading
# apostrophe to force text storage, then the cell style is reset to the
# workbook's default "Normal" so no stray number-format/quote-prefix
# styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.984.38'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').Value = '3.485.04'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'602.24"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.73%  '
$ws.Range('D6').Value = "'148.37"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.49%  '
$ws.Range('D7').Value = '3.482.08'
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.479"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D11').Value = "'7.70"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.90%  '
$ws.Range('E12').Value = '  -3.43%  '
$ws.Range('D13').Value = "'0.0000213"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('D14').Value = '4.072.42'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = "'31.18"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.14%  '
$ws.Range('D16').Value = '3.480.02'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '66.897.92'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('E19').Value = '  -5.04%  '
$ws.Range('D20').Value = "'10.21"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('D21').Value = "'15.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.93%  '
$ws.Range('D22').Value = "'435.01"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = '  -5.86%  '
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D26').Value = '3.620.46'
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('E27').Value = '  -9.98%  '
$ws.Range('D28').Value = "'9.83"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.82%  '
$ws.Range('D29').Value = "'8.37"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.43%  '
$ws.Range('D30').Value = "'2.49"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').Value = "'1.60"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.92%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = "'0.168"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').Value = "'0.998"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').Value = "'25.36"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('D35').Value = '3.473.32'
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = "'5.93"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.28%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'1.80"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.32%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = "'7.90"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.18%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = "'173.05"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.60%  '
$ws.Range('D42').Value = "'0.0885"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.61%  '
$ws.Range('D43').Value = "'2.09"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -12.60%  '
$ws.Range('D44').Value = "'5.40"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.55%  '
$ws.Range('D45').Value = "'0.900"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = "'29.01"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.04%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = "'46.47"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('E48').Value = '  -6.97%  '
$ws.Range('E49').Value = '  -4.55%  '
$ws.Range('D50').Value = "'2.42"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.39%  '
$ws.Range('D51').Value = "'0.971"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.98%  '
